$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '42.861.65'
$cell.Style = $origStyle
$ws.Range('E2').Value = '  -0.59%  '
$cell = $ws.Range('D3')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.211.85'
$cell.Style = $origStyle
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  -0.11%  '
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '255.75'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  +1.51%  '
$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.617'
$cell.Style = $origStyle
$ws.Range('E6').Value = '  +0.37%  '
$cell = $ws.Range('D7')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '77.23'
$cell.Style = $origStyle
$ws.Range('E7').Value = '  +2.60%  '
$ws.Range('E8').Value = '  -0.07%  '
$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.592'
$cell.Style = $origStyle
$ws.Range('E9').Value = '  -1.51%  '
$cell = $ws.Range('D10')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '42.69'
$cell.Style = $origStyle
$ws.Range('E10').Value = '  +2.92%  '
$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0909'
$cell.Style = $origStyle
$ws.Range('E11').Value = '  -2.19%  '
$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '6.99'
$cell.Style = $origStyle
$ws.Range('E12').Value = '  +1.13%  '
$cell = $ws.Range('D13')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.102'
$cell.Style = $origStyle
$ws.Range('E13').Value = '  +0.66%  '
$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.540.64'
$cell.Style = $origStyle
$ws.Range('E14').Value = '  -1.41%  '
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '14.41'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  -1.75%  '
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.209.58'
$cell.Style = $origStyle
$ws.Range('E16').Value = '  -1.79%  '
$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.783'
$cell.Style = $origStyle
$ws.Range('E17').Value = '  -1.22%  '
$cell = $ws.Range('D18')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '42.802.91'
$cell.Style = $origStyle
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('E19').Value = '  -0.97%  '
$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '71.03'
$cell.Style = $origStyle
$ws.Range('E20').Value = '  -0.20%  '
$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.95'
$cell.Style = $origStyle
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('E22').Value = '  +6.26%  '
$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '229.89'
$cell.Style = $origStyle
$ws.Range('E23').Value = '  +0.40%  '
$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '9.21'
$cell.Style = $origStyle
$ws.Range('E24').Value = '  -4.54%  '
$ws.Range('E25').Value = '  -0.13%  '
$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '42.60'
$cell.Style = $origStyle
$ws.Range('E26').Value = '  +8.05%  '
$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '10.70'
$cell.Style = $origStyle
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('E28').Value = '  -2.93%  '
$ws.Range('E29').Value = '  -2.87%  '
$ws.Range('E30').Value = '  -1.25%  '
$cell = $ws.Range('D31')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '172.33'
$cell.Style = $origStyle
$ws.Range('E31').Value = '  +0.24%  '
$cell = $ws.Range('D32')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0879'
$cell.Style = $origStyle
$ws.Range('E32').Value = '  +9.56%  '
$cell = $ws.Range('D33')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '20.34'
$cell.Style = $origStyle
$ws.Range('E33').Value = '  +0.48%  '
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('E35').Value = '  -0.23%  '
$cell = $ws.Range('D36')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0357'
$cell.Style = $origStyle
$ws.Range('E36').Value = '  +7.58%  '
$ws.Range('E37').Value = '  -3.27%  '
$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.39'
$cell.Style = $origStyle
$ws.Range('E38').Value = '  -2.47%  '
$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '13.11'
$cell.Style = $origStyle
$ws.Range('E39').Value = '  +1.08%  '
$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.89'
$cell.Style = $origStyle
$ws.Range('E40').Value = '  +19.00%  '
$ws.Range('E41').Value = '  -0.21%  '
$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '60.95'
$cell.Style = $origStyle
$ws.Range('E42').Value = '  +2.25%  '
$ws.Range('E43').Value = '  -1.90%  '
$ws.Range('E44').Value = '  -2.74%  '
$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '103.03'
$cell.Style = $origStyle
$ws.Range('E45').Value = '  -0.76%  '
$cell = $ws.Range('D46')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '8.44'
$cell.Style = $origStyle
$ws.Range('E46').Value = '  -2.61%  '
$cell = $ws.Range('D47')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0973'
$cell.Style = $origStyle
$ws.Range('E47').Value = '  -1.61%  '
$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.464'
$cell.Style = $origStyle
$ws.Range('E48').Value = '  -3.92%  '
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.50'
$cell.Style = $origStyle
$ws.Range('E50').Value = '  +23.81%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.13'
$cell.Style = $origStyle
$ws.Range('E51').Value = '  -1.49%  '
